$p = $ppt.ActivePresentation
$s = $p.Slides.Item(13)
$shape = $s.Shapes.Item(4)
$textRange = $shape.TextFrame.TextRange
$find = $textRange.Replace("17:30 - 15:30 UTC", "17:30 - 19:30 UTC")
